# Fix issues resulting from extra files: the sheet previously mixed data
# rows from several different tickers (WIX, NOC, CDNS, ...). Restore every
# row to the EBAY-only open/close/high/low price, shares_outstanding (H)
# and fixed_ticker (I) values so the sheet - and therefore the shared
# string table, once the now-unused ticker strings are dropped on save -
# only ever describes EBAY.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = 21.54948957371399
$ws.Cells.Item(2, 5).Value = 21.89906311035156
$ws.Cells.Item(2, 6).Value = 22.6057258426151
$ws.Cells.Item(2, 7).Value = 20.82403088014538
$ws.Cells.Item(2, 8).Value = 457000000
$ws.Cells.Item(2, 9).Value = "EBAY"

$ws.Cells.Item(3, 4).Value = 22.92522175871082
$ws.Cells.Item(3, 5).Value = 25.11404609680176
$ws.Cells.Item(3, 6).Value = 26.21256189626325
$ws.Cells.Item(3, 7).Value = 22.48543658985835
$ws.Cells.Item(3, 8).Value = 457000000
$ws.Cells.Item(3, 9).Value = "EBAY"

$ws.Cells.Item(4, 4).Value = 21.8810156255608
$ws.Cells.Item(4, 5).Value = 24.91756439208984
$ws.Cells.Item(4, 6).Value = 25.51594325488637
$ws.Cells.Item(4, 7).Value = 21.09508601791939
$ws.Cells.Item(4, 8).Value = 457000000
$ws.Cells.Item(4, 9).Value = "EBAY"

$ws.Cells.Item(5, 4).Value = 24.21200809222352
$ws.Cells.Item(5, 5).Value = 20.95218276977539
$ws.Cells.Item(5, 6).Value = 24.33704189077607
$ws.Cells.Item(5, 7).Value = 20.00549464475374
$ws.Cells.Item(5, 8).Value = 457000000
$ws.Cells.Item(5, 9).Value = "EBAY"

$ws.Cells.Item(6, 4).Value = 21.76490721995367
$ws.Cells.Item(6, 5).Value = 21.8184928894043
$ws.Cells.Item(6, 6).Value = 23.40821414294883
$ws.Cells.Item(6, 7).Value = 20.89859676263694
$ws.Cells.Item(6, 8).Value = 457000000
$ws.Cells.Item(6, 9).Value = "EBAY"

$ws.Cells.Item(7, 4).Value = 20.92539553995367
$ws.Cells.Item(7, 5).Value = 27.82907867431641
$ws.Cells.Item(7, 6).Value = 28.39173430253645
$ws.Cells.Item(7, 7).Value = 20.80929124053004
$ws.Cells.Item(7, 8).Value = 457000000
$ws.Cells.Item(7, 9).Value = "EBAY"

$ws.Cells.Item(8, 4).Value = 27.95411415063723
$ws.Cells.Item(8, 5).Value = 25.46235847473145
$ws.Cells.Item(8, 6).Value = 29.49918143060584
$ws.Cells.Item(8, 7).Value = 25.39984155756459
$ws.Cells.Item(8, 8).Value = 457000000
$ws.Cells.Item(8, 9).Value = "EBAY"

$ws.Cells.Item(9, 4).Value = 26.64124890126123
$ws.Cells.Item(9, 5).Value = 28.42745399475098
$ws.Cells.Item(9, 6).Value = 29.33841729780272
$ws.Cells.Item(9, 7).Value = 26.35545635885613
$ws.Cells.Item(9, 8).Value = 457000000
$ws.Cells.Item(9, 9).Value = "EBAY"

$ws.Cells.Item(10, 4).Value = 29.80283517986716
$ws.Cells.Item(10, 5).Value = 29.83856010437012
$ws.Cells.Item(10, 6).Value = 31.02638829222093
$ws.Cells.Item(10, 7).Value = 28.48104363550475
$ws.Cells.Item(10, 8).Value = 457000000
$ws.Cells.Item(10, 9).Value = "EBAY"

$ws.Cells.Item(11, 4).Value = 31.39254777731885
$ws.Cells.Item(11, 5).Value = 31.91054534912109
$ws.Cells.Item(11, 6).Value = 33.47347440247721
$ws.Cells.Item(11, 7).Value = 30.32082431649166
$ws.Cells.Item(11, 8).Value = 457000000
$ws.Cells.Item(11, 9).Value = "EBAY"

$ws.Cells.Item(12, 4).Value = 34.39337703365663
$ws.Cells.Item(12, 5).Value = 33.61637878417969
$ws.Cells.Item(12, 6).Value = 35.08106639263399
$ws.Cells.Item(12, 7).Value = 32.32137943044092

$ws.Cells.Item(13, 4).Value = 34.08972672120942
$ws.Cells.Item(13, 5).Value = 36.24210739135742
$ws.Cells.Item(13, 6).Value = 36.9208639104165
$ws.Cells.Item(13, 7).Value = 33.38417821270161
$ws.Cells.Item(13, 8).Value = 457000000
$ws.Cells.Item(13, 9).Value = "EBAY"

$ws.Cells.Item(14, 4).Value = 35.75983612103065
$ws.Cells.Item(14, 5).Value = 33.83073425292969
$ws.Cells.Item(14, 6).Value = 38.40341953159184
$ws.Cells.Item(14, 7).Value = 33.62532102240467
$ws.Cells.Item(14, 8).Value = 457000000
$ws.Cells.Item(14, 9).Value = "EBAY"

$ws.Cells.Item(15, 4).Value = 32.07131238127867
$ws.Cells.Item(15, 5).Value = 29.87428092956543
$ws.Cells.Item(15, 6).Value = 34.24155184674922
$ws.Cells.Item(15, 7).Value = 29.76710957767879
$ws.Cells.Item(15, 8).Value = 457000000
$ws.Cells.Item(15, 9).Value = "EBAY"

$ws.Cells.Item(16, 4).Value = 29.56169706052044
$ws.Cells.Item(16, 5).Value = 25.92677116394043
$ws.Cells.Item(16, 6).Value = 30.90135106761804
$ws.Cells.Item(16, 7).Value = 23.54218696316838
$ws.Cells.Item(16, 8).Value = 457000000
$ws.Cells.Item(16, 9).Value = "EBAY"

$ws.Cells.Item(17, 4).Value = 24.6942867763111
$ws.Cells.Item(17, 5).Value = 30.05290412902832
$ws.Cells.Item(17, 6).Value = 31.03531563529274
$ws.Cells.Item(17, 7).Value = 24.61390740529557
$ws.Cells.Item(17, 8).Value = 457000000
$ws.Cells.Item(17, 9).Value = "EBAY"

$ws.Cells.Item(18, 4).Value = 33.53590778380505
$ws.Cells.Item(18, 5).Value = 34.73714065551758
$ws.Cells.Item(18, 6).Value = 34.78196208856122
$ws.Cells.Item(18, 7).Value = 31.76991528070554
$ws.Cells.Item(18, 8).Value = 457000000
$ws.Cells.Item(18, 9).Value = "EBAY"

$ws.Cells.Item(19, 4).Value = 35.97029922006883
$ws.Cells.Item(19, 5).Value = 37.06821441650391
$ws.Cells.Item(19, 6).Value = 37.79716089826816
$ws.Cells.Item(19, 7).Value = 35.0343686977228
$ws.Cells.Item(19, 8).Value = 457000000
$ws.Cells.Item(19, 9).Value = "EBAY"

$ws.Cells.Item(20, 4).Value = 35.20998997479518
$ws.Cells.Item(20, 5).Value = 31.83257484436035
$ws.Cells.Item(20, 6).Value = 35.64345441035537
$ws.Cells.Item(20, 7).Value = 31.72420959668758
$ws.Cells.Item(20, 8).Value = 457000000
$ws.Cells.Item(20, 9).Value = "EBAY"

$ws.Cells.Item(21, 4).Value = 33.00902008865313
$ws.Cells.Item(21, 5).Value = 30.42523384094238
$ws.Cells.Item(21, 6).Value = 33.09061352060056
$ws.Cells.Item(21, 7).Value = 30.39803372471017
$ws.Cells.Item(21, 8).Value = 457000000
$ws.Cells.Item(21, 9).Value = "EBAY"

$ws.Cells.Item(22, 4).Value = 26.72142689697695
$ws.Cells.Item(22, 5).Value = 36.27520370483398
$ws.Cells.Item(22, 6).Value = 36.75790048611601
$ws.Cells.Item(22, 7).Value = 25.77424529191476
$ws.Cells.Item(22, 8).Value = 457000000
$ws.Cells.Item(22, 9).Value = "EBAY"

$ws.Cells.Item(23, 4).Value = 48.17996005491695
$ws.Cells.Item(23, 5).Value = 50.52908706665039
$ws.Cells.Item(23, 6).Value = 55.81234069846441
$ws.Cells.Item(23, 7).Value = 46.90942094975978
$ws.Cells.Item(23, 8).Value = 457000000
$ws.Cells.Item(23, 9).Value = "EBAY"

$ws.Cells.Item(24, 4).Value = 48.01072817977752
$ws.Cells.Item(24, 5).Value = 43.66528701782227
$ws.Cells.Item(24, 6).Value = 53.48378740616621
$ws.Cells.Item(24, 7).Value = 43.07856108010362
$ws.Cells.Item(24, 8).Value = 457000000
$ws.Cells.Item(24, 9).Value = "EBAY"

$ws.Cells.Item(25, 4).Value = 46.42272655099118
$ws.Cells.Item(25, 5).Value = 51.96807098388672
$ws.Cells.Item(25, 6).Value = 55.68336138539236
$ws.Cells.Item(25, 7).Value = 46.39513888972353
$ws.Cells.Item(25, 8).Value = 457000000
$ws.Cells.Item(25, 9).Value = "EBAY"

$ws.Cells.Item(26, 4).Value = 56.64398924933778
$ws.Cells.Item(26, 5).Value = 51.46853637695312
$ws.Cells.Item(26, 6).Value = 60.06661372870931
$ws.Cells.Item(26, 7).Value = 50.69360220056591
$ws.Cells.Item(26, 8).Value = 457000000
$ws.Cells.Item(26, 9).Value = "EBAY"

$ws.Cells.Item(27, 4).Value = 64.96807996691787
$ws.Cells.Item(27, 5).Value = 63.11740112304688
$ws.Cells.Item(27, 6).Value = 68.59540880655794
$ws.Cells.Item(27, 7).Value = 62.25683517825568
$ws.Cells.Item(27, 8).Value = 457000000
$ws.Cells.Item(27, 9).Value = "EBAY"

$ws.Cells.Item(28, 4).Value = 64.64619744867976
$ws.Cells.Item(28, 5).Value = 71.15719604492188
$ws.Cells.Item(28, 6).Value = 75.3030869224207
$ws.Cells.Item(28, 7).Value = 63.32916034644617
$ws.Cells.Item(28, 8).Value = 457000000
$ws.Cells.Item(28, 9).Value = "EBAY"

$ws.Cells.Item(29, 4).Value = 61.78948060006648
$ws.Cells.Item(29, 5).Value = 55.85694885253906
$ws.Cells.Item(29, 6).Value = 62.75654033438191
$ws.Cells.Item(29, 7).Value = 51.98871346242819
$ws.Cells.Item(29, 8).Value = 457000000
$ws.Cells.Item(29, 9).Value = "EBAY"

$ws.Cells.Item(30, 4).Value = 54.03423084150439
$ws.Cells.Item(30, 5).Value = 48.47861099243164
$ws.Cells.Item(30, 6).Value = 54.50108959248782
$ws.Cells.Item(30, 7).Value = 48.10512612875466
$ws.Cells.Item(30, 8).Value = 457000000
$ws.Cells.Item(30, 9).Value = "EBAY"

$ws.Cells.Item(31, 4).Value = 39.33738914528377
$ws.Cells.Item(31, 5).Value = 45.61223983764648
$ws.Cells.Item(31, 6).Value = 46.64397730671171
$ws.Cells.Item(31, 7).Value = 38.66207157086663
$ws.Cells.Item(31, 8).Value = 457000000
$ws.Cells.Item(31, 9).Value = "EBAY"

$ws.Cells.Item(32, 4).Value = 34.82000028348973
$ws.Cells.Item(32, 5).Value = 37.55356979370117
$ws.Cells.Item(32, 6).Value = 38.29823230385974
$ws.Cells.Item(32, 7).Value = 33.85853797843131
$ws.Cells.Item(32, 8).Value = 457000000
$ws.Cells.Item(32, 9).Value = "EBAY"

$ws.Cells.Item(33, 4).Value = 39.86539494905932
$ws.Cells.Item(33, 5).Value = 46.89488983154297
$ws.Cells.Item(33, 6).Value = 47.26436414289514
$ws.Cells.Item(33, 7).Value = 39.35381346537068
$ws.Cells.Item(33, 8).Value = 457000000
$ws.Cells.Item(33, 9).Value = "EBAY"

$ws.Cells.Item(34, 4).Value = 42.06404489647473
$ws.Cells.Item(34, 5).Value = 44.23632049560547
$ws.Cells.Item(34, 6).Value = 44.46497917211179
$ws.Cells.Item(34, 7).Value = 40.48247511852683
$ws.Cells.Item(34, 8).Value = 457000000
$ws.Cells.Item(34, 9).Value = "EBAY"

$ws.Cells.Item(35, 4).Value = 42.79522748991629
$ws.Cells.Item(35, 5).Value = 42.6514892578125
$ws.Cells.Item(35, 6).Value = 47.41396873814946
$ws.Cells.Item(35, 7).Value = 41.49201352428094
$ws.Cells.Item(35, 8).Value = 457000000
$ws.Cells.Item(35, 9).Value = "EBAY"

$ws.Cells.Item(36, 4).Value = 42.31292812566448
$ws.Cells.Item(36, 5).Value = 37.80314636230469
$ws.Cells.Item(36, 6).Value = 42.43820086241089
$ws.Cells.Item(36, 7).Value = 36.55043002269387
$ws.Cells.Item(36, 8).Value = 457000000
$ws.Cells.Item(36, 9).Value = "EBAY"

$ws.Cells.Item(37, 4).Value = 41.94458948476228
$ws.Cells.Item(37, 5).Value = 39.82118225097656
$ws.Cells.Item(37, 6).Value = 42.856009227001
$ws.Cells.Item(37, 7).Value = 38.93885281242229
$ws.Cells.Item(37, 8).Value = 457000000
$ws.Cells.Item(37, 9).Value = "EBAY"

$ws.Cells.Item(38, 4).Value = 51.32195006269598
$ws.Cells.Item(38, 5).Value = 50.23994827270508
$ws.Cells.Item(38, 6).Value = 51.58513627031154
$ws.Cells.Item(38, 7).Value = 48.16367466808671
$ws.Cells.Item(38, 8).Value = 457000000
$ws.Cells.Item(38, 9).Value = "EBAY"

$ws.Cells.Item(39, 4).Value = 52.75868645567719
$ws.Cells.Item(39, 5).Value = 54.4830207824707
$ws.Cells.Item(39, 6).Value = 54.72795433841608
$ws.Cells.Item(39, 7).Value = 50.97557046738604
$ws.Cells.Item(39, 8).Value = 457000000
$ws.Cells.Item(39, 9).Value = "EBAY"

$ws.Cells.Item(40, 4).Value = 64.12514845329184
$ws.Cells.Item(40, 5).Value = 56.60532760620117
$ws.Cells.Item(40, 6).Value = 66.73346368381313
$ws.Cells.Item(40, 7).Value = 55.44389324926963
$ws.Cells.Item(40, 8).Value = 457000000
$ws.Cells.Item(40, 9).Value = "EBAY"

$ws.Cells.Item(41, 4).Value = 61.20290912737551
$ws.Cells.Item(41, 5).Value = 66.69852447509766
$ws.Cells.Item(41, 6).Value = 70.69173102028779
$ws.Cells.Item(41, 7).Value = 60.40229206603858
$ws.Cells.Item(41, 8).Value = 457000000
$ws.Cells.Item(41, 9).Value = "EBAY"

$ws.Cells.Item(42, 4).Value = 67.12008738709555
$ws.Cells.Item(42, 5).Value = 67.67611694335938
$ws.Cells.Item(42, 6).Value = 68.48036410524296
$ws.Cells.Item(42, 7).Value = 58.29320056204837
$ws.Cells.Item(42, 8).Value = 457000000
$ws.Cells.Item(42, 9).Value = "EBAY"

$ws.Cells.Item(43, 4).Value = 74.55628091884633
$ws.Cells.Item(43, 5).Value = 91.46327972412109
$ws.Cells.Item(43, 6).Value = 92.50003061949134
$ws.Cells.Item(43, 7).Value = 74.42668990900725
$ws.Cells.Item(43, 8).Value = 457000000
$ws.Cells.Item(43, 9).Value = "EBAY"
